$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -1.517586829436955
$ws.Range("F2").Value = 0.07587934147184777
$ws.Range("G2").Value = 0.1490118507824147
$ws.Range("H2").Value = 0.1490118507824147
$ws.Range("J2").Value = 0.0
$ws.Range("R2").Value = 1.517586829436955
$ws.Range("S2").Value = -0.0
$ws.Range("U2").Value = 0.0
$ws.Range("V2").Value = 0.0
$ws.Range("C3").Value = 52.41206585281522
$ws.Range("D3").Value = 0.7891180289236067
$ws.Range("F3").Value = 0.03945590144618034
$ws.Range("G3").Value = -0.07682853129600235
$ws.Range("J3").Value = 0.0
$ws.Range("P3").Value = 0.7891180289236067
$ws.Range("S3").Value = -0.0
$ws.Range("U3").Value = 0.0
$ws.Range("V3").Value = 0.0
$ws.Range("C4").Value = 56.35765599743326
$ws.Range("D4").Value = 5.473934116369112
$ws.Range("F4").Value = 0.2736967058184556
$ws.Range("G4").Value = -0.5317926994052593
$ws.Range("J4").Value = 0.0
$ws.Range("P4").Value = 5.473934116369112
$ws.Range("S4").Value = -0.0
$ws.Range("U4").Value = 0.0
$ws.Range("V4").Value = 0.0
$ws.Range("C5").Value = 83.72732657927881
$ws.Range("D5").Value = 1.254534684144237
$ws.Range("F5").Value = 0.06272673420721186
$ws.Range("G5").Value = -0.121263322569382
$ws.Range("J5").Value = 0.0
$ws.Range("P5").Value = 1.254534684144237
$ws.Range("S5").Value = -0.0
$ws.Range("U5").Value = 0.0
$ws.Range("V5").Value = 0.0
$ws.Range("C6").Value = 90.0
$ws.Range("D6").Value = 0.0
$ws.Range("F6").Value = 0.0
$ws.Range("G6").Value = 0.0
$ws.Range("H6").Value = 0.0
$ws.Range("J6").Value = 0.0
$ws.Range("R6").Value = 0.0
$ws.Range("S6").Value = -0.0
$ws.Range("U6").Value = 0.0
$ws.Range("V6").Value = 0.0
$ws.Range("C7").Value = 90.0
$ws.Range("D7").Value = 0.0
$ws.Range("F7").Value = 0.0
$ws.Range("G7").Value = 0.0
$ws.Range("P7").Value = 0.0
$ws.Range("V7").Value = 0.0
$ws.Range("C8").Value = 90.0
$ws.Range("G8").Value = 0.0
$ws.Range("J8").Value = 0.0
$ws.Range("S8").Value = -0.0
$ws.Range("U8").Value = 0.0
$ws.Range("V8").Value = 0.0
$ws.Range("C9").Value = 90.0
$ws.Range("D9").Value = -4.195753223810646
$ws.Range("F9").Value = 0.2097876611905323
$ws.Range("G9").Value = 0.4959799885866565
$ws.Range("H9").Value = 0.4959799885866565
$ws.Range("J9").Value = 0.0
$ws.Range("R9").Value = 4.195753223810646
$ws.Range("S9").Value = -0.0
$ws.Range("U9").Value = 0.0
$ws.Range("V9").Value = 0.0
$ws.Range("C10").Value = 69.02123388094677
$ws.Range("D10").Value = -4.015976791473788
$ws.Range("F10").Value = 0.2007988395736894
$ws.Range("G10").Value = 0.5823166347636993
$ws.Range("H10").Value = 0.5823166347636993
$ws.Range("J10").Value = 0.0
$ws.Range("R10").Value = 4.015976791473788
$ws.Range("S10").Value = -0.0
$ws.Range("U10").Value = 0.0
$ws.Range("V10").Value = 0.0
$ws.Range("C11").Value = 48.94134992357782
$ws.Range("D11").Value = -3.826771154896456
$ws.Range("F11").Value = 0.1913385577448228
$ws.Range("G11").Value = 0.4952607228666992
$ws.Range("H11").Value = 0.4952607228666992
$ws.Range("J11").Value = 0.0
$ws.Range("R11").Value = 3.826771154896456
$ws.Range("S11").Value = -0.0
$ws.Range("U11").Value = 0.0
$ws.Range("V11").Value = 0.0
$ws.Range("C12").Value = 29.80749414909555
$ws.Range("D12").Value = -1.961498829819109
$ws.Range("F12").Value = 0.09807494149095547
$ws.Range("G12").Value = 0.1977583120223626
$ws.Range("H12").Value = 0.1977583120223626
$ws.Range("J12").Value = 0.0
$ws.Range("R12").Value = 1.961498829819109
$ws.Range("S12").Value = -0.0
$ws.Range("U12").Value = 0.0
$ws.Range("V12").Value = 0.0
$ws.Range("C13").Value = 20.0
$ws.Range("D13").Value = 0.01468583641022033
$ws.Range("F13").Value = 0.0007342918205110166
$ws.Range("G13").Value = -0.001382230922929938
$ws.Range("H13").Value = 0.0
$ws.Range("J13").Value = 0.0
$ws.Range("P13").Value = 0.01468583641022033
$ws.Range("R13").Value = 0.0
$ws.Range("S13").Value = -0.0
$ws.Range("U13").Value = 0.0
$ws.Range("V13").Value = 0.0
$ws.Range("C14").Value = 20.0734291820511
$ws.Range("D14").Value = 9.570397286250321
$ws.Range("F14").Value = 0.478519864312516
$ws.Range("G14").Value = -0.7905148158442764
$ws.Range("P14").Value = 9.570397286250321
$ws.Range("V14").Value = 0.0
$ws.Range("C15").Value = 67.9254156133027
$ws.Range("D15").Value = 4.414916877339459
$ws.Range("F15").Value = 0.2207458438669729
$ws.Range("G15").Value = -0.3779168847002576
$ws.Range("J15").Value = 0.0
$ws.Range("P15").Value = 4.414916877339459
$ws.Range("S15").Value = -0.0
$ws.Range("U15").Value = 0.0
$ws.Range("V15").Value = 0.0
$ws.Range("C16").Value = 90.0
$ws.Range("D16").Value = 0.0
$ws.Range("F16").Value = 0.0
$ws.Range("G16").Value = 0.0
$ws.Range("H16").Value = 0.0
$ws.Range("J16").Value = 0.0
$ws.Range("R16").Value = 0.0
$ws.Range("S16").Value = -0.0
$ws.Range("U16").Value = 0.0
$ws.Range("V16").Value = 0.0
$ws.Range("C17").Value = 90.0
$ws.Range("D17").Value = 0.0
$ws.Range("F17").Value = 0.0
$ws.Range("G17").Value = 0.0
$ws.Range("P17").Value = 0.0
$ws.Range("V17").Value = 0.0
$ws.Range("C18").Value = 90.0
$ws.Range("G18").Value = 0.0
$ws.Range("J18").Value = 0.0
$ws.Range("S18").Value = -0.0
$ws.Range("U18").Value = 0.0
$ws.Range("V18").Value = 0.0
$ws.Range("C19").Value = 90.0
$ws.Range("G19").Value = 0.0
$ws.Range("J19").Value = 0.0
$ws.Range("S19").Value = -0.0
$ws.Range("U19").Value = 0.0
$ws.Range("V19").Value = 0.0
$ws.Range("C20").Value = 90.0
$ws.Range("G20").Value = 0.0
$ws.Range("J20").Value = 0.0
$ws.Range("S20").Value = -0.0
$ws.Range("U20").Value = 0.0
$ws.Range("V20").Value = 0.0
$ws.Range("C21").Value = 90.0
$ws.Range("D21").Value = -4.195753223810646
$ws.Range("F21").Value = 0.2097876611905323
$ws.Range("G21").Value = 0.6335587367954076
$ws.Range("H21").Value = 0.6335587367954076
$ws.Range("J21").Value = 0.0
$ws.Range("R21").Value = 4.195753223810646
$ws.Range("S21").Value = -0.0
$ws.Range("U21").Value = 0.0
$ws.Range("V21").Value = 0.0
$ws.Range("C22").Value = 69.02123388094677
$ws.Range("D22").Value = -4.015976791473788
$ws.Range("F22").Value = 0.2007988395736894
$ws.Range("G22").Value = 0.6017941222023471
$ws.Range("H22").Value = 0.6017941222023471
$ws.Range("J22").Value = 0.0
$ws.Range("R22").Value = 4.015976791473788
$ws.Range("S22").Value = -0.0
$ws.Range("U22").Value = 0.0
$ws.Range("V22").Value = 0.0
$ws.Range("C23").Value = 48.94134992357782
$ws.Range("D23").Value = -3.826771154896456
$ws.Range("F23").Value = 0.1913385577448228
$ws.Range("G23").Value = 0.5120219805251458
$ws.Range("H23").Value = 0.5120219805251458
$ws.Range("J23").Value = 0.0
$ws.Range("R23").Value = 3.826771154896456
$ws.Range("S23").Value = -0.0
$ws.Range("U23").Value = 0.0
$ws.Range("V23").Value = 0.0
$ws.Range("C24").Value = 29.80749414909555
$ws.Range("D24").Value = -1.961498829819109
$ws.Range("F24").Value = 0.09807494149095547
$ws.Range("G24").Value = 0.2037016534767145
$ws.Range("H24").Value = 0.2037016534767145
$ws.Range("J24").Value = 0.0
$ws.Range("R24").Value = 1.961498829819109
$ws.Range("S24").Value = -0.0
$ws.Range("U24").Value = 0.0
$ws.Range("V24").Value = 0.0
$ws.Range("C25").Value = 20.0
$ws.Range("D25").Value = 2.430912548402899
$ws.Range("F25").Value = 0.1215456274201449
$ws.Range("G25").Value = -0.2333676046466783
$ws.Range("J25").Value = 0.0
$ws.Range("P25").Value = 2.430912548402899
$ws.Range("S25").Value = -0.0
$ws.Range("U25").Value = 0.0
$ws.Range("V25").Value = 0.0
$ws.Range("C26").Value = 32.15456274201449
$ws.Range("D26").Value = -2.430912548402898
$ws.Range("F26").Value = 0.1215456274201449
$ws.Range("G26").Value = 0.2502624468580784
$ws.Range("H26").Value = 0.2502624468580784
$ws.Range("J26").Value = 0.0
$ws.Range("R26").Value = 2.430912548402898
$ws.Range("S26").Value = -0.0
$ws.Range("U26").Value = 0.0
$ws.Range("V26").Value = 0.0
$ws.Range("C27").Value = 20.0
$ws.Range("G27").Value = 0.0
$ws.Range("V27").Value = 0.0
